# The sheet had a small "Name / Email / Orcid" table in A1:C2, with B2
# carrying a mailto: hyperlink (styled with the built-in "Hyperlink" style).
# The edit clears all the cell contents (the example data) while leaving
# B2's Hyperlink formatting/style in place (now just an empty, styled
# cell), and removes the now-unused hyperlink itself.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the hyperlink on B2 (ClearContents alone does not drop hyperlinks).
foreach ($hl in $ws.Hyperlinks) {
    $hl.Delete()
}

# Clear all the example values (Name/Email/Orcid header + Marie Monfils row)
# but keep cell formatting/styles intact - this leaves B2 with its Hyperlink
# style (s="1") but no value, matching a plain "select all, press Delete".
$ws.UsedRange.ClearContents()

# Drop the now fully-empty, un-styled row 1 so it doesn't linger in the
# saved sheetData with stale row-height metadata.
$ws.Rows(1).AutoFit()

# Leave the sheet selection the way it ended up after clearing the example
# data out - a block spanning the old table plus a little extra room.
[void]$ws.Range("A1:D5").Select()
